# Update "想去人数" (want-to-go count) figures across the four sheets to the
# freshly regenerated values, as published by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 67
$ws.Range("F5").Value  = 1943
$ws.Range("F7").Value  = 7652
$ws.Range("F16").Value = 3688
$ws.Range("F17").Value = 5897
$ws.Range("F20").Value = 1026
$ws.Range("F22").Value = 387
$ws.Range("F23").Value = 6005
$ws.Range("F29").Value = 1881
$ws.Range("F31").Value = 273
$ws.Range("F35").Value = 313
$ws.Range("F38").Value = 1839
$ws.Range("F42").Value = 1074
$ws.Range("F44").Value = 549
$ws.Range("F47").Value = 75

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 14
$ws.Range("F11").Value = 341
$ws.Range("F12").Value = 392

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value  = 839
$ws.Range("F11").Value = 1106
$ws.Range("F12").Value = 1491

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 1943
$ws.Range("F8").Value  = 7652
$ws.Range("F14").Value = 1106
$ws.Range("F18").Value = 1491
$ws.Range("F19").Value = 3688
$ws.Range("F20").Value = 341
$ws.Range("F21").Value = 392
$ws.Range("F23").Value = 1026
$ws.Range("F25").Value = 387
$ws.Range("F26").Value = 6005
$ws.Range("F29").Value = 1881
$ws.Range("F31").Value = 273
$ws.Range("F35").Value = 313
$ws.Range("F38").Value = 1839
$ws.Range("F42").Value = 1074
$ws.Range("F44").Value = 549
$ws.Range("F46").Value = 75
